$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text changes: the single "strikePrice" column is split into a
#    CE strike-price column (H) and a PE strike-price column (new, P).
#    Set P1 first so the existing shared string "strikePrice" is freed/renamed
#    to "strikePrice_PE" at its current index, then H1 gets a new shared
#    string "strikePrice_CE" appended.
# ---------------------------------------------------------------------------
$ws.Range("P1").Value = "strikePrice_PE"
$ws.Range("H1").Value = "strikePrice_CE"

# ---------------------------------------------------------------------------
# 2. Column widths (new spacer column A, resized columns E/H/O/P, and newly
#    explicit widths for I/J/K which previously used the sheet default).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 5.5
$ws.Columns.Item(5).ColumnWidth = 11
$ws.Columns.Item(8).ColumnWidth = 11.16666666666667
$ws.Columns.Item(9).ColumnWidth = 6
$ws.Columns.Item(10).ColumnWidth = 16
$ws.Columns.Item(11).ColumnWidth = 13.5
$ws.Columns.Item(15).ColumnWidth = 17.5
$ws.Columns.Item(16).ColumnWidth = 11.66666666666667

# ---------------------------------------------------------------------------
# 3. Header row formatting: center-align + shrink-to-fit across A1:P1.
#    Build each combined style once (on a cell that already carries the
#    right base number format) and fan it out with copy/paste-special so we
#    do not litter the style table with one-property-at-a-time variants.
# ---------------------------------------------------------------------------
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").ShrinkToFit = $true
$ws.Range("A1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("I1:M1").PasteSpecial(-4122)

$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").ShrinkToFit = $true
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("N1:O1").PasteSpecial(-4122)
$ws.Range("P1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. View state: clear the scrolled-down top-left cell and move the active
#    selection to H6.
# ---------------------------------------------------------------------------
$ws.Range("H6").Select()

Write-Output "edit applied"
